$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

# Sheet 1 (展览) column F ("想去人数") updates
$ws1.Cells.Item(5, 6).Value = 11
$ws1.Cells.Item(6, 6).Value = 199
$ws1.Cells.Item(7, 6).Value = 4597
$ws1.Cells.Item(10, 6).Value = 104
$ws1.Cells.Item(12, 6).Value = 88
$ws1.Cells.Item(13, 6).Value = 706
$ws1.Cells.Item(14, 6).Value = 190
$ws1.Cells.Item(15, 6).Value = 1001
$ws1.Cells.Item(16, 6).Value = 79
$ws1.Cells.Item(17, 6).Value = 240
$ws1.Cells.Item(18, 6).Value = 156
$ws1.Cells.Item(21, 6).Value = 95
$ws1.Cells.Item(22, 6).Value = 3560
$ws1.Cells.Item(23, 6).Value = 5911
$ws1.Cells.Item(25, 6).Value = 33
$ws1.Cells.Item(29, 6).Value = 3370
$ws1.Cells.Item(30, 6).Value = 367
$ws1.Cells.Item(32, 6).Value = 2477
$ws1.Cells.Item(35, 6).Value = 126
$ws1.Cells.Item(36, 6).Value = 217
$ws1.Cells.Item(37, 6).Value = 267
$ws1.Cells.Item(39, 6).Value = 130
$ws1.Cells.Item(40, 6).Value = 1020
$ws1.Cells.Item(43, 6).Value = 29
$ws1.Cells.Item(45, 6).Value = 55
$ws1.Cells.Item(47, 6).Value = 65
$ws1.Cells.Item(48, 6).Value = 551

# Sheet 4 (全部类型) column F ("想去人数") updates
$ws4.Cells.Item(5, 6).Value = 11
$ws4.Cells.Item(6, 6).Value = 199
$ws4.Cells.Item(7, 6).Value = 4597
$ws4.Cells.Item(10, 6).Value = 104
$ws4.Cells.Item(13, 6).Value = 88
$ws4.Cells.Item(14, 6).Value = 706
$ws4.Cells.Item(15, 6).Value = 190
$ws4.Cells.Item(16, 6).Value = 1001
$ws4.Cells.Item(17, 6).Value = 79
$ws4.Cells.Item(18, 6).Value = 240
$ws4.Cells.Item(19, 6).Value = 156
$ws4.Cells.Item(22, 6).Value = 95
$ws4.Cells.Item(23, 6).Value = 3560
$ws4.Cells.Item(24, 6).Value = 5912
$ws4.Cells.Item(26, 6).Value = 33
$ws4.Cells.Item(30, 6).Value = 3370
$ws4.Cells.Item(31, 6).Value = 367
$ws4.Cells.Item(33, 6).Value = 2477
$ws4.Cells.Item(36, 6).Value = 126
$ws4.Cells.Item(37, 6).Value = 217
$ws4.Cells.Item(38, 6).Value = 267
$ws4.Cells.Item(40, 6).Value = 130
$ws4.Cells.Item(41, 6).Value = 1020
$ws4.Cells.Item(44, 6).Value = 29
$ws4.Cells.Item(46, 6).Value = 55
$ws4.Cells.Item(48, 6).Value = 65
$ws4.Cells.Item(49, 6).Value = 551

Write-Output "applied updates"
